# Apply the "output generated at 456a3b4" update to 北京-漫展信息.xlsx
#
# Summary of the change:
#  - Sheet 1 (展览)      : bump "想去人数" (F column) for a batch of rows (re-scrape).
#  - Sheet 2 (演出)      : one listing (2024-07-22 石川绫子小提琴动漫音乐会, row 5) was
#                          removed from the source feed. All following rows shift up by
#                          one; the row-index column A keeps its original 0-based
#                          sequence (0,1,2,...,22), so only columns B:I move.
#  - Sheet 3 (本地生活)  : bump F column for a few rows.
#  - Sheet 4 (全部类型)  : bump F column for a batch of rows (mirrors sheet 1 and the
#                          sheet-2 "地点" row whose count rose independently of the
#                          deletion).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 - 展览 : "想去人数" (column F) increases
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    3  = 126
    6  = 241
    7  = 12893
    8  = 46
    10 = 228
    11 = 2955
    13 = 6293
    14 = 62
    16 = 3338
    17 = 23
    18 = 158
    22 = 56
    24 = 3567
    25 = 77
    27 = 2699
    28 = 386
    29 = 1856
    30 = 97
    31 = 202
    32 = 6503
    33 = 14
    34 = 162
    35 = 214
    36 = 1946
    37 = 1288
    39 = 1012
    40 = 13
    41 = 203
    42 = 214
    45 = 122
    46 = 1178
    47 = 1730
    48 = 147
    49 = 1163
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2 - 演出 : remove the row-5 listing ("石川绫子小提琴动漫音乐会", id=83973).
# Excel's native row delete shifts every column (including A) up by one; the
# source index column A must stay fixed at its original sequence, so after
# deleting we restore A5:A23 = 4..22.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(5).Delete()
for ($row = 5; $row -le 23; $row++) {
    $ws2.Cells.Item($row, 1).Value = $row - 1
}

# ---------------------------------------------------------------------------
# Sheet 3 - 本地生活 : "想去人数" (column F) increases
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$sheet3Updates = @{
    2 = 419
    3 = 576
    4 = 8
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 4 - 全部类型 : "想去人数" (column F) increases
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    3  = 126
    6  = 419
    7  = 576
    8  = 241
    10 = 12893
    13 = 228
    14 = 2956
    16 = 6293
    17 = 62
    20 = 56
    22 = 121
    23 = 3568
    24 = 77
    27 = 2700
    28 = 1856
    29 = 97
    30 = 202
    31 = 6503
    33 = 162
    34 = 215
    35 = 1946
    37 = 1288
    39 = 1012
    40 = 203
    41 = 214
    44 = 122
    45 = 1178
    47 = 1730
    48 = 147
    49 = 1163
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
